$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Row 55: the XRP sell that was "IN PROGRESS" got CANCELled, with a
# finalized-date note typed into column I (entered as text, not a date serial).
# ---------------------------------------------------------------------------
$ws.Range("H55").Value = "CANCEL"
$ws.Range("I55").Value = " 2017-05-12 20:35:22"

# ---------------------------------------------------------------------------
# Row 56: brand-new trade entered right underneath it.
# Pull the per-column formatting (A/D/I) from row 55 first, since those
# three columns carry direct (non-default) cell styles in this sheet.
# ---------------------------------------------------------------------------
$ws.Range("A55").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("D55").Copy()
$ws.Range("D56").PasteSpecial(-4122)
$ws.Range("I55").Copy()
$ws.Range("I56").PasteSpecial(-4122)

# D56 holds a number-looking string ("0.1947"); force Text so it is stored
# as a label instead of being auto-parsed into a numeric cell.
$ws.Range("D56").NumberFormat = "@"

$ws.Range("A56").Value = " 2017-05-12 20:35:22"
$ws.Range("B56").Value = "            Buy"
$ws.Range("C56").Value = "        XRP"
$ws.Range("D56").Value = "        0.1947`n"
$ws.Range("E56").Value = "         0.185USDT"
$ws.Range("F56").Value = "         179 XRP"
$ws.Range("G56").Value = " XRP/USDT0000006"
$ws.Range("H56").Value = "IN PROGRESS"
$ws.Range("I56").Value = " "
$ws.Range("K56").Value = "     "

# Re-apply D55's number format/alignment over D56 (the "@" tweak above was
# only needed to land the value as text) so the visible formatting matches
# the rest of the column again.
$ws.Range("D55").Copy()
$ws.Range("D56").PasteSpecial(-4122)

# The embedded newline triggers wrap-driven autofit; pin the row back to the
# same 14.25pt height used throughout the table.
$ws.Rows.Item(56).RowHeight = 14.25

# ---------------------------------------------------------------------------
# View state: selection moved down onto the newly-entered row.
# ---------------------------------------------------------------------------
$ws.Range("D59").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
